$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top; this pushes the 63 existing data rows
# down to rows 2-64 (dimension grows from A1:E63 to A1:E64).
$ws.Rows.Item(1).Insert()

# Write the new header row (becomes shared strings X1, X2, X3, X4, Y).
$ws.Range("A1").Value = "X1"
$ws.Range("B1").Value = "X2"
$ws.Range("C1").Value = "X3"
$ws.Range("D1").Value = "X4"
$ws.Range("E1").Value = "Y"

# The inserted row copied formatting down from the original row 1 (now row
# 2, style index 1 for A-D, style index 2 for E). The header row should use
# style index 1 across all five columns, so copy the A:D formatting (style 1)
# over the whole header range A1:E1.
$ws.Range("A2:D2").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)

# Update the saved selection/view: no more frozen/scrolled topLeftCell, and
# the active cell moves to H8.
[void]$ws.Range("H8").Select()
